$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (the "Förändrad" date column) rows 2-97 all need to be updated
# from serial date 45172 (2023-09-03) to 45175 (2023-09-06).
$ws.Range("C2:C97").Value = 45175
